$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22
$ws.Range("A22").Formula = '="2.20.2020"'
$ws.Range("B22").Value = '5:00 - 8:00'
$ws.Range("C22").Value = 'Class'
$ws.Range("D22").Value = 'Architectual Understanding'
$ws.Range("E22").Value = 'Being able to understand that Architectual design of a program'
$ws.Range("F22").Value = 'It was an extremely tough process for me. Realizing that I wasn''t trying to identify how the program should be architect, but the fact that we had to identify how the program is architected as is made the activity much more challening.'
$ws.Range("G22").Value = 'After seeing that the homework was to come up with the architectual model of our open-source program, I dreaded it. It felt daunting to me that I struggled so much with PacMan, and now to do that with our huge project, it seemed near impossible to do. I am not sure how to even approach things.'

# Row 23
$ws.Range("A23").Formula = '="2.21.2020"'
$ws.Range("B23").Value = '2:00 - 7:30'
$ws.Range("C23").Value = 'Nic, Rafael, Chris'
$ws.Range("D23").Value = 'Architectual Model of Open Source'
$ws.Range("E23").Value = 'We were able to generate the architectual model and find issues/pull requests that were interesting'
$ws.Range("F23").Value = 'The assignment was actual not as hard as I had envisioned, unlike how I felt the day before. Once we abstracted most of the functionality of the program, it was a lot easier to see the overall architectual model. '
$ws.Range("G23").Value = 'I personally enjoyed finding some of the issues that people ran into, and it was cool '

# Row 24
$ws.Range("A24").Formula = '="2.27.2020"'
$ws.Range("B24").Value = '5:00 - 8:00 pm'
$ws.Range("C24").Value = 'Class'
$ws.Range("D24").Value = 'Architectual Patterns'
$ws.Range("E24").Value = 'Learned about forms of architectual practices / patterns that are well built common coding practices that has knowledge built on them for common coding problems'
$ws.Range("F24").Value = 'The class tied very closely with Professor Malek''s Architecture course, and drew a lot of knowledge we knew from that class to this topic.'
$ws.Range("G24").Value = 'Actually utilizing something we learned in one class and applying it to better our knowledge of how other systems worked felt very useful. While it''s hard for me to truly grasp what some architectual patterns are, it probably will become more clear to me as I continue to grow and become more exposed to future coding challenges.'

# Row 25
$ws.Range("A25").Formula = '="3.1.2020"'
$ws.Range("B25").Value = '2:00 - 7:00 pm '
$ws.Range("C25").Value = 'Nic, Rafael, Chris'
$ws.Range("D25").Value = 'Identifying Architectual Patterns in Code, Solve issue'
$ws.Range("E25").Value = 'was able to identify and solve an issue within FreeCol, and find a couple of the architectual patterns. '
$ws.Range("F25").Value = 'Started by trying to find the easiest bug to fix. Given our current understanding of the game, it would be both extremely hard to identify and test a bug that is embedded in the gaming logic. We decide to look for bugs that were centered around the map editor, as that is something we had access to as soon as we open the game. Not only did we get to tackle an issue others seem to be having, but getting a thumbs up for our bug fix felt very reassuring and ultimately felt great!'
$ws.Range("G25").Value = 'Given how small our bug was, it seems crazy to think how long it might take others to solve certain issues / bugs that are exponentially larger in scale. It also feels a little daunting to try and tackle a second issue. I can imagine how hard it must be for others with much more complicated projects probably struggle a lot more than us.'

# Convert column-A date formulas into literal text values, preserving cell format
$dateRange = $ws.Range("A22:A25")
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Row heights to match the wrapped multi-line content
$ws.Rows(22).RowHeight = 115.2
$ws.Rows(23).RowHeight = 109.2
$ws.Rows(24).RowHeight = 129.6
$ws.Rows(25).RowHeight = 234

# Update the view: scroll position and active selection
$win = $wb.Windows.Item(1)
$win.ScrollRow = 23
$win.ScrollColumn = 1
$ws.Range("G25").Select()
